$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix C11 to hold the "Changes were made..." text (was previously "Title was choosen")
$ws.Range("C11").Value = "Changes were made to the document accordingly based on the review"

# Add new row 12 entries
$ws.Range("A12").Value = "17/12/2019"
$ws.Range("B12").Value = "dfd and changes in document was made"
$ws.Range("C12").Value = "dfd"

# Update selection to C12, matching the saved view state
$ws.Activate()
$ws.Range("C12").Select()
